$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Build the four new cell-format combinations in a scratch area (column Z)
#    so that they become new entries in cellXfs, then copy/paste-special the
#    formatting onto the real destination cells.  The scratch cells are
#    cleared at the end so they do not remain part of the sheet.
# ---------------------------------------------------------------------------

# Style A: Courier New, 7pt, black, vertical-center               (future xf "G4")
$ws.Range("Z1").Value = "stage"
$ws.Range("Z1").Font.Size = 7
$ws.Range("Z1").Font.Color = 0
$ws.Range("Z1").Font.Name = "Courier New"
$ws.Range("Z1").VerticalAlignment = -4108

# Style B: Segoe UI, 7pt, dark-gray, no special alignment          (future xf "G5")
$ws.Range("Z2").Value = "stage"
$ws.Range("Z2").Font.Size = 7
$ws.Range("Z2").Font.Color = 2171169
$ws.Range("Z2").Font.Name = "Segoe UI"

# Style C: Segoe UI, 7pt, dark-gray, left/center, wrap             (future xf "G6")
$ws.Range("Z3").Value = "stage"
$ws.Range("Z3").Font.Size = 7
$ws.Range("Z3").Font.Color = 2171169
$ws.Range("Z3").Font.Name = "Segoe UI"
$ws.Range("Z3").HorizontalAlignment = -4131
$ws.Range("Z3").VerticalAlignment = -4108
$ws.Range("Z3").WrapText = $true
$ws.Range("Z3").ReadingOrder = 1

# Style D: Courier New, 7pt, black, vertical-center, wrap          (future xf "G7")
$ws.Range("Z4").Value = "stage"
$ws.Range("Z4").Font.Size = 7
$ws.Range("Z4").Font.Color = 0
$ws.Range("Z4").Font.Name = "Courier New"
$ws.Range("Z4").VerticalAlignment = -4108
$ws.Range("Z4").WrapText = $true

# ---------------------------------------------------------------------------
# 2. Column E got wider (stored width needs to end up as exactly 41)
# ---------------------------------------------------------------------------
$ws.Columns("E").ColumnWidth = 40.14

# ---------------------------------------------------------------------------
# 3. Row 6 - edit note (success / not-found)
# ---------------------------------------------------------------------------
$ws.Range("B6").Value = "http://localhost:3030/editnote/6839bcc70c7906ec01e05903"
$ws.Range("B6").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("B6"), "http://localhost:3030/editnote/6839bcc70c7906ec01e05903") | Out-Null
$ws.Range("B6").Style = "Hyperlink"

$ws.Range("C6").Value = @'
{
    "Message":"testmessage d",
    "userId":"6839bab50c7906ec01e04e78"
}
'@
$ws.Range("C6").WrapText = $true

$ws.Range("D6").Value = @'
{
    "status": "Note updated",
    "note": {
        "_id": "6839d3040cf40c3a33c3eaf7",
        "userId": "6839bab50c7906ec01e04e78",
        "Message": "testmessage d",
        "notesDate": "2025-05-30T15:47:16.662Z",
        "__v": 0
    }
}
'@
$ws.Range("D6").WrapText = $true

$ws.Range("E6").Value = @'
{
    "status": "Note not found"
}
'@
$ws.Range("E6").WrapText = $true

$ws.Rows("6:6").RowHeight = 259.2

# ---------------------------------------------------------------------------
# 4. Row 7 - search notes
# ---------------------------------------------------------------------------
$ws.Range("B7").Value = "http://localhost:3030/searchnotes"
$ws.Range("Z2").Copy()
$ws.Range("B7").PasteSpecial(-4122)

$ws.Range("C7").Value = @'
{
    "userId": "6839bab50c7906ec01e04e78",
    "notesDate": "2025-05-29"
}
'@
$ws.Range("C7").WrapText = $true

$ws.Range("D7").Value = @'
{
    "status": "success",
    "items": []
}
'@
$ws.Range("D7").WrapText = $true

$ws.Range("E7").Value = @'
{
    "status": "Error fetching notes",}
'@
$ws.Range("E7").WrapText = $true

$ws.Rows("7:7").RowHeight = 86.4

# ---------------------------------------------------------------------------
# 5. Row 8 - delete note
# ---------------------------------------------------------------------------
$ws.Range("B8").Value = "http://localhost:3030/deletenote/6839d3040cf40c3a33c3eaf7"
$ws.Range("Z3").Copy()
$ws.Range("B8").PasteSpecial(-4122)

$ws.Range("C8").Value = '  "_id": "683991463baaa70086a79dce"}'
$ws.Range("Z1").Copy()
$ws.Range("C8").PasteSpecial(-4122)

$ws.Range("D8").Value = @'
{
    "status": "Note deleted"
}
'@
$ws.Range("D8").WrapText = $true

$ws.Range("E8").Value = '{    "status": "Error",}'
$ws.Range("E8").WrapText = $true

$ws.Rows("8:8").RowHeight = 57.6

# ---------------------------------------------------------------------------
# 6. Row 9 - edit note (second example, reusing row 6's JSON payloads)
# ---------------------------------------------------------------------------
$ws.Range("B9").Value = "http://localhost:3030/editnote/6839d3040cf40c3a33c3eaf7"
$ws.Range("Z2").Copy()
$ws.Range("B9").PasteSpecial(-4122)

$ws.Range("C9").Value = $ws.Range("C6").Value()
$ws.Range("Z4").Copy()
$ws.Range("C9").PasteSpecial(-4122)

$ws.Range("D9").Value = $ws.Range("D6").Value()
$ws.Range("D9").WrapText = $true

$ws.Range("E9").Value = $ws.Range("E6").Value()
$ws.Range("Z4").Copy()
$ws.Range("E9").PasteSpecial(-4122)

$ws.Rows("9:9").RowHeight = 259.2

# ---------------------------------------------------------------------------
# 7. Row 10 - two empty, styled cells
# ---------------------------------------------------------------------------
$ws.Range("Z1").Copy()
$ws.Range("C10").PasteSpecial(-4122)
$ws.Range("Z1").Copy()
$ws.Range("E10").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 8. Clean up scratch cells
# ---------------------------------------------------------------------------
$ws.Range("Z1:Z4").Clear()

# ---------------------------------------------------------------------------
# 9. View state: mirror the author's final selection / scroll position
# ---------------------------------------------------------------------------
$ws.Range("B9").Select()
